$wb = $excel.ActiveWorkbook

# --- L1_PES ---
$ws = $wb.Worksheets.Item("L1_PES")
$ws.Range("D2").Value = 0.0602
$ws.Range("D3").Value = 0.0602
$ws.Range("D4").Value = 0.0602
$ws.Range("I4").Value = 0.0254
$ws.Range("N4").Value = 0.007
$ws.Range("O4").Value = 0.0086
$ws.Range("P4").Value = 0.0086
$ws.Range("Q4").Value = 0.0258
$ws.Range("R4").Value = 0.0086
$ws.Range("S4").Value = 0.0086
$ws.Range("T4").Value = 0.0086
$ws.Range("U4").Value = 0.0258
$ws.Range("V4").Value = 0.1032
$ws.Range("D5").Value = 0.592233009708738
$ws.Range("D6").Value = 0.592233009708738
$ws.Range("D7").Value = 0.592233009708738
$ws.Range("N7").Value = 0.8261
$ws.Range("O7").Value = 0.592233009708738
$ws.Range("P7").Value = 0.592233009708738
$ws.Range("Q7").Value = 0.592233009708738
$ws.Range("R7").Value = 0.592233009708738
$ws.Range("S7").Value = 0.592233009708738
$ws.Range("T7").Value = 0.592233009708738
$ws.Range("U7").Value = 0.592233009708738
$ws.Range("V7").Value = 0.592233009708738

# --- PES APAC ---
$ws = $wb.Worksheets.Item("PES APAC")
$ws.Range("D2").Value = 0.0389
$ws.Range("D3").Value = 0.0389
$ws.Range("D4").Value = 0.0389
$ws.Range("N4").Value = 0.0072
$ws.Range("O4").Value = 0.00555833333333333
$ws.Range("P4").Value = 0.00555833333333333
$ws.Range("Q4").Value = 0.016675
$ws.Range("R4").Value = 0.00555833333333333
$ws.Range("S4").Value = 0.00555833333333333
$ws.Range("T4").Value = 0.00555833333333333
$ws.Range("U4").Value = 0.016675
$ws.Range("V4").Value = 0.0667
$ws.Range("D5").Value = 0.2
$ws.Range("D6").Value = 0.2
$ws.Range("D7").Value = 0.2
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 0.2
$ws.Range("P7").Value = 0.2
$ws.Range("Q7").Value = 0.2
$ws.Range("R7").Value = 0.2
$ws.Range("S7").Value = 0.2
$ws.Range("T7").Value = 0.2
$ws.Range("U7").Value = 0.2
$ws.Range("V7").Value = 0.2

# --- PES EMEA ---
$ws = $wb.Worksheets.Item("PES EMEA")
$ws.Range("D2").Value = 0.0136
$ws.Range("D3").Value = 0.0136
$ws.Range("D4").Value = 0.0136
$ws.Range("I4").Value = 0.0091
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0.00194166666666667
$ws.Range("P4").Value = 0.00194166666666667
$ws.Range("Q4").Value = 0.005825
$ws.Range("R4").Value = 0.00194166666666667
$ws.Range("S4").Value = 0.00194166666666667
$ws.Range("T4").Value = 0.00194166666666667
$ws.Range("U4").Value = 0.005825
$ws.Range("V4").Value = 0.0233
$ws.Range("D5").Value = 0.862068965517241
$ws.Range("D6").Value = 0.862068965517241
$ws.Range("D7").Value = 0.862068965517241
$ws.Range("M7").Value = 0.8696
$ws.Range("N7").Value = 0.8
$ws.Range("O7").Value = 0.862068965517241
$ws.Range("P7").Value = 0.862068965517241
$ws.Range("Q7").Value = 0.862068965517241
$ws.Range("R7").Value = 0.862068965517241
$ws.Range("S7").Value = 0.862068965517241
$ws.Range("T7").Value = 0.862068965517241
$ws.Range("U7").Value = 0.862068965517241
$ws.Range("V7").Value = 0.862068965517241

# --- PES NA Motors Solutions ---
$ws = $wb.Worksheets.Item("PES NA Motors Solutions")
$ws.Range("D2").Value = 0.0762
$ws.Range("D3").Value = 0.0762
$ws.Range("D4").Value = 0.0762
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0.0108833333333333
$ws.Range("P4").Value = 0.0108833333333333
$ws.Range("Q4").Value = 0.03265
$ws.Range("R4").Value = 0.0108833333333333
$ws.Range("S4").Value = 0.0108833333333333
$ws.Range("T4").Value = 0.0108833333333333
$ws.Range("U4").Value = 0.03265
$ws.Range("V4").Value = 0.1306
$ws.Range("N7").ClearContents()

# --- PES NA Motors and Drives ---
$ws = $wb.Worksheets.Item("PES NA Motors and Drives")
$ws.Range("D2").Value = 0.0648
$ws.Range("D3").Value = 0.0648
$ws.Range("D4").Value = 0.0648
$ws.Range("F4").Value = 0.012
$ws.Range("G4").Value = 0.0201
$ws.Range("I4").Value = 0.0361
$ws.Range("K4").Value = 0.0102
$ws.Range("M4").Value = 0.0203
$ws.Range("N4").Value = 0.0082
$ws.Range("O4").Value = 0.00925833333333333
$ws.Range("P4").Value = 0.00925833333333333
$ws.Range("Q4").Value = 0.027775
$ws.Range("R4").Value = 0.00925833333333333
$ws.Range("S4").Value = 0.00925833333333333
$ws.Range("T4").Value = 0.00925833333333333
$ws.Range("U4").Value = 0.027775
$ws.Range("V4").Value = 0.1111
$ws.Range("D5").Value = 0.473684210526316
$ws.Range("D6").Value = 0.473684210526316
$ws.Range("D7").Value = 0.473684210526316
$ws.Range("J7").Value = 0.1667
$ws.Range("M7").Value = 0.3077
$ws.Range("N7").Value = 0.9167
$ws.Range("O7").Value = 0.473684210526316
$ws.Range("P7").Value = 0.473684210526316
$ws.Range("Q7").Value = 0.473684210526316
$ws.Range("R7").Value = 0.473684210526316
$ws.Range("S7").Value = 0.473684210526316
$ws.Range("T7").Value = 0.473684210526316
$ws.Range("U7").Value = 0.473684210526316
$ws.Range("V7").Value = 0.473684210526316

# --- PES NA Sales ---
$ws = $wb.Worksheets.Item("PES NA Sales")
$ws.Range("D2").Value = 0.1624
$ws.Range("D3").Value = 0.1624
$ws.Range("D4").Value = 0.1624
$ws.Range("F4").Value = 0.0288
$ws.Range("G4").Value = 0.0294
$ws.Range("I4").Value = 0.0575
$ws.Range("J4").Value = 0.028
$ws.Range("K4").Value = 0.0374
$ws.Range("L4").Value = 0.0288
$ws.Range("M4").Value = 0.0943
$ws.Range("N4").Value = 0.0098
$ws.Range("O4").Value = 0.0232
$ws.Range("P4").Value = 0.0232
$ws.Range("Q4").Value = 0.0696
$ws.Range("R4").Value = 0.0232
$ws.Range("S4").Value = 0.0232
$ws.Range("T4").Value = 0.0232
$ws.Range("U4").Value = 0.0696
$ws.Range("V4").Value = 0.2784
$ws.Range("D5").Value = 0.625
$ws.Range("D6").Value = 0.625
$ws.Range("D7").Value = 0.625
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 0.625
$ws.Range("P7").Value = 0.625
$ws.Range("Q7").Value = 0.625
$ws.Range("R7").Value = 0.625
$ws.Range("S7").Value = 0.625
$ws.Range("T7").Value = 0.625
$ws.Range("U7").Value = 0.625
$ws.Range("V7").Value = 0.625

# --- PES NA Strategy & Mktg ---
$ws = $wb.Worksheets.Item("PES NA Strategy & Mktg")
$ws.Range("D2").Value = 0.0643
$ws.Range("D3").Value = 0.0643
$ws.Range("D4").Value = 0.0643
$ws.Range("F4").Value = 0.011
$ws.Range("I4").Value = 0.0106
$ws.Range("J4").Value = 0.0171
$ws.Range("L4").Value = 0.0081
$ws.Range("M4").Value = 0.0251
$ws.Range("N4").Value = 0.0248
$ws.Range("O4").Value = 0.00918333333333333
$ws.Range("P4").Value = 0.00918333333333333
$ws.Range("Q4").Value = 0.02755
$ws.Range("R4").Value = 0.00918333333333333
$ws.Range("S4").Value = 0.00918333333333333
$ws.Range("T4").Value = 0.00918333333333333
$ws.Range("U4").Value = 0.02755
$ws.Range("V4").Value = 0.1102
$ws.Range("D5").Value = 0.8
$ws.Range("D6").Value = 0.8
$ws.Range("D7").Value = 0.8
$ws.Range("J7").Value = 1
$ws.Range("M7").Value = 0.75
$ws.Range("N7").ClearContents()
$ws.Range("O7").Value = 0.8
$ws.Range("P7").Value = 0.8
$ws.Range("Q7").Value = 0.8
$ws.Range("R7").Value = 0.8
$ws.Range("S7").Value = 0.8
$ws.Range("T7").Value = 0.8
$ws.Range("U7").Value = 0.8
$ws.Range("V7").Value = 0.8

# --- PES Segment Functions ---
$ws = $wb.Worksheets.Item("PES Segment Functions")
$ws.Range("D2").Value = 0.1097
$ws.Range("D3").Value = 0.1097
$ws.Range("D4").Value = 0.1097
$ws.Range("H4").Value = 0.0357
$ws.Range("I4").Value = 0.0349
$ws.Range("J4").Value = 0.037
$ws.Range("K4").Value = 0.0189
$ws.Range("L4").Value = 0.0192
$ws.Range("M4").Value = 0.0755
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0.015675
$ws.Range("P4").Value = 0.015675
$ws.Range("Q4").Value = 0.047025
$ws.Range("R4").Value = 0.015675
$ws.Range("S4").Value = 0.015675
$ws.Range("T4").Value = 0.015675
$ws.Range("U4").Value = 0.047025
$ws.Range("V4").Value = 0.1881
$ws.Range("D5").Value = 0.6
$ws.Range("D6").Value = 0.6
$ws.Range("D7").Value = 0.6
$ws.Range("J7").ClearContents()
$ws.Range("M7").Value = 0.5
$ws.Range("N7").ClearContents()
$ws.Range("O7").Value = 0.6
$ws.Range("P7").Value = 0.6
$ws.Range("Q7").Value = 0.6
$ws.Range("R7").Value = 0.6
$ws.Range("S7").Value = 0.6
$ws.Range("T7").Value = 0.6
$ws.Range("U7").Value = 0.6
$ws.Range("V7").Value = 0.6
